$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.978.89'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').Value = '2.466.81'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '516.98'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -3.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.30'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -4.15%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.99%  '
$ws.Range('D9').Value = '2.469.12'
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0990'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('E13').Value = '  -2.47%  '
$ws.Range('D14').Value = '2.905.57'
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '57.955.09'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.31'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -3.38%  '
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = '2.466.77'
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.73'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.17'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '318.86'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  -4.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.18'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.409'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -3.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.32'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('D29').Value = '0.0₃0747'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('E30').Value = '  -4.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '164.78'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.23'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -6.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.16'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -8.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.97'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('E39').Value = '  -4.91%  '
$ws.Range('E40').Value = '  -2.96%  '
$ws.Range('E41').Value = '  -4.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '271.75'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -4.88%  '
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '126.56'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -3.58%  '
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0210'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.91'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').Value = '1.723.58'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('E51').Value = '  -1.83%  '
